$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2301
$ws.Range("J17").Value = 2301
$ws.Range("L17").Value = 6903
$ws.Range("N17").Value = -7239
# Row 43
$ws.Range("H43").Value = 4273.636
$ws.Range("I43").Value = 5960.6924
$ws.Range("J43").Value = 1836.7778
$ws.Range("K43").Value = 5960.6924
$ws.Range("L43").Value = 1836.7778
$ws.Range("M43").Value = -5891.6924
$ws.Range("N43").Value = -1974.7778
# Row 62
$ws.Range("H62").Value = 8932.23
$ws.Range("I62").Value = 9601.727999999999
$ws.Range("J62").Value = 5250
$ws.Range("K62").Value = 9601.727999999999
$ws.Range("L62").Value = 5250
$ws.Range("M62").Value = -8977.727999999999
$ws.Range("N62").Value = -6498
# Row 65
$ws.Range("H65").Value = 8932.23
$ws.Range("I65").Value = 9601.727999999999
$ws.Range("J65").Value = 5250
$ws.Range("K65").Value = 48008.64
$ws.Range("L65").Value = 26250
$ws.Range("M65").Value = -44888.64
$ws.Range("N65").Value = -32490
# Row 127
$ws.Range("H127").Value = 1399

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4937.915
$ws.Range("I32").Value = 4063.849
$ws.Range("K32").Value = 4063.849
$ws.Range("M32").Value = -3776.849
# Row 61
$ws.Range("H61").Value = 3964.8076
$ws.Range("I61").Value = 2178
$ws.Range("J61").Value = 6401.364
$ws.Range("K61").Value = 2178
$ws.Range("L61").Value = 6401.364
$ws.Range("M61").Value = -1966
$ws.Range("N61").Value = -6825.364
# Row 88
$ws.Range("H88").Value = 8500
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 8500
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 8500
$ws.Range("N88").Value = -9312
$ws.Range("M88").ClearContents()
# Row 91
$ws.Range("H91").Value = 8500
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 8500
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 8500
$ws.Range("N91").Value = -11308
$ws.Range("M91").ClearContents()
# Row 122
$ws.Range("H122").Value = 2745.6567
$ws.Range("I122").Value = 2397.2166
$ws.Range("K122").Value = 7191.649800000001
$ws.Range("M122").Value = -4741.649800000001
# Row 132
$ws.Range("H132").Value = 2279.1904
$ws.Range("I132").Value = 1612.25
$ws.Range("J132").Value = 3613.0715
$ws.Range("K132").Value = 4836.75
$ws.Range("L132").Value = 10839.2145
$ws.Range("M132").Value = -2306.75
$ws.Range("N132").Value = -15899.2145
# Row 136
$ws.Range("H136").Value = 3964.8076
$ws.Range("I136").Value = 2178
$ws.Range("J136").Value = 6401.364
$ws.Range("K136").Value = 6534
$ws.Range("L136").Value = 19204.092
$ws.Range("M136").Value = -3984
$ws.Range("N136").Value = -24304.092

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3249.0715
$ws.Range("I86").Value = 2767.1052
$ws.Range("J86").Value = 4266.5557
$ws.Range("K86").Value = 2767.1052
$ws.Range("L86").Value = 4266.5557
$ws.Range("M86").Value = -1644.1052
$ws.Range("N86").Value = -6512.5557
# Row 89
$ws.Range("H89").Value = 3249.0715
$ws.Range("I89").Value = 2767.1052
$ws.Range("J89").Value = 4266.5557
$ws.Range("K89").Value = 13835.526
$ws.Range("L89").Value = 21332.7785
$ws.Range("M89").Value = -8219.526
$ws.Range("N89").Value = -32564.7785
# Row 105
$ws.Range("H105").Value = 8668777
$ws.Range("I105").Value = 436580.56
$ws.Range("J105").Value = 35717420
$ws.Range("K105").Value = 436580.56
$ws.Range("L105").Value = 35717420
$ws.Range("M105").Value = -434833.56
$ws.Range("N105").Value = -35720914
# Row 107
$ws.Range("H107").Value = 3078462
$ws.Range("I107").Value = 4050027
$ws.Range("J107").Value = 1838.6666
$ws.Range("K107").Value = 4050027
$ws.Range("L107").Value = 1838.6666
$ws.Range("M107").Value = -4048107
$ws.Range("N107").Value = -5678.6666
# Row 134
$ws.Range("H134").Value = 2992.6
$ws.Range("I134").Value = 2244.75
$ws.Range("K134").Value = 6734.25
$ws.Range("M134").Value = -4199.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2490
$ws.Range("I16").Value = 2490
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2490
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2203
$ws.Range("N16").ClearContents()
# Row 31
$ws.Range("H31").Value = 3777.2354
$ws.Range("I31").Value = 3098.7
$ws.Range("J31").Value = 8866.25
$ws.Range("K31").Value = 3098.7
$ws.Range("L31").Value = 8866.25
$ws.Range("M31").Value = -2803.7
$ws.Range("N31").Value = -9456.25
# Row 34
$ws.Range("H34").Value = 3777.2354
$ws.Range("I34").Value = 3098.7
$ws.Range("J34").Value = 8866.25
$ws.Range("K34").Value = 3098.7
$ws.Range("L34").Value = 8866.25
$ws.Range("M34").Value = -2896.7
$ws.Range("N34").Value = -9270.25
# Row 113
$ws.Range("H113").Value = 2490
$ws.Range("I113").Value = 2490
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2490
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -320
$ws.Range("N113").ClearContents()
# Row 132
$ws.Range("H132").Value = 11907583
$ws.Range("I132").Value = 13891235
$ws.Range("J132").Value = 5666.5
$ws.Range("K132").Value = 41673705
$ws.Range("L132").Value = 16999.5
$ws.Range("M132").Value = -41671175
$ws.Range("N132").Value = -22059.5
# Row 134
$ws.Range("H134").Value = 3116.158
$ws.Range("I134").Value = 2762.9375
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8288.8125
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -5753.8125
$ws.Range("N134").Value = -20070

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4770.8887
$ws.Range("I3").Value = 4347.375
$ws.Range("K3").Value = 13042.125
$ws.Range("M3").Value = -12930.125
# Row 34
$ws.Range("H34").Value = 1992.5385
$ws.Range("J34").Value = 2316.6667
$ws.Range("L34").Value = 6950.000100000001
$ws.Range("N34").Value = -7118.000100000001
# Row 39
$ws.Range("H39").Value = 8424.75
$ws.Range("J39").Value = 9599.714
$ws.Range("L39").Value = 28799.142
$ws.Range("N39").Value = -29387.142
# Row 55
$ws.Range("H55").Value = 9443.429
$ws.Range("J55").Value = 10867.333
$ws.Range("L55").Value = 32601.999
$ws.Range("N55").Value = -32955.999
# Row 56
$ws.Range("H56").Value = 7030.25
$ws.Range("I56").Value = 7030.25
$ws.Range("K56").Value = 7030.25
$ws.Range("M56").Value = -6500.25
# Row 92
$ws.Range("H92").Value = 661.75
$ws.Range("I92").Value = 999
$ws.Range("K92").Value = 2997
$ws.Range("M92").Value = -1749
# Row 113
$ws.Range("H113").Value = 6609.6
$ws.Range("I113").Value = 565
$ws.Range("J113").Value = 7676.294
$ws.Range("K113").Value = 1695
$ws.Range("L113").Value = 23028.882
$ws.Range("M113").Value = 475
$ws.Range("N113").Value = -27368.882
# Row 125
$ws.Range("H125").Value = 3010
$ws.Range("I125").Value = 2015
$ws.Range("K125").Value = 6045
$ws.Range("M125").Value = -1125
# Row 126
$ws.Range("H126").Value = 3300
$ws.Range("I126").Value = 3300
$ws.Range("K126").Value = 9900
$ws.Range("M126").Value = -4960

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 68134.2
$ws.Range("J57").Value = 68134.2
$ws.Range("L57").Value = 68134.2
$ws.Range("N57").Value = -69774.2
# Row 80
$ws.Range("H80").Value = 50001696
$ws.Range("I80").Value = 111112450
$ws.Range("J80").Value = 1986.8182
$ws.Range("K80").Value = 111112450
$ws.Range("L80").Value = 1986.8182
$ws.Range("M80").Value = -111111452
$ws.Range("N80").Value = -3982.8182
# Row 83
$ws.Range("H83").Value = 50001696
$ws.Range("I83").Value = 111112450
$ws.Range("J83").Value = 1986.8182
$ws.Range("K83").Value = 555562250
$ws.Range("L83").Value = 9934.091
$ws.Range("M83").Value = -555557258
$ws.Range("N83").Value = -19918.091
# Row 102
$ws.Range("H102").Value = 4731.34
$ws.Range("I102").Value = 787.4091
$ws.Range("K102").Value = 787.4091
$ws.Range("M102").Value = 834.5909
# Row 132
$ws.Range("H132").Value = 7857
$ws.Range("I132").Value = 12999.667
$ws.Range("K132").Value = 38999.001
$ws.Range("M132").Value = -36469.001
# Row 134
$ws.Range("H134").Value = 54813.43
$ws.Range("J134").Value = 54813.43
$ws.Range("L134").Value = 164440.29
$ws.Range("N134").Value = -169510.29

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1570.7693
$ws.Range("I82").Value = 2243.875
$ws.Range("K82").Value = 2243.875
$ws.Range("M82").Value = -1882.875
# Row 85
$ws.Range("H85").Value = 1570.7693
$ws.Range("I85").Value = 2243.875
$ws.Range("K85").Value = 2243.875
$ws.Range("M85").Value = -995.875
# Row 136
$ws.Range("H136").Value = 4017.5625
$ws.Range("I136").Value = 3711.5
$ws.Range("K136").Value = 11134.5
$ws.Range("M136").Value = -8584.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 24900
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 24900
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 113
$ws.Range("H113").Value = 798
$ws.Range("I113").Value = 852.8889
$ws.Range("J113").Value = 727.4286
$ws.Range("K113").Value = 2558.6667
$ws.Range("L113").Value = 2182.2858
$ws.Range("M113").Value = -388.6667000000002
$ws.Range("N113").Value = -6522.2858
# Row 126
$ws.Range("H126").Value = 2566.8667
$ws.Range("I126").Value = 2375.25
$ws.Range("K126").Value = 7125.75
$ws.Range("M126").Value = -4655.75
# Row 132
$ws.Range("H132").Value = 6805137
$ws.Range("I132").Value = 8335912
$ws.Range("K132").Value = 25007736
$ws.Range("M132").Value = -25005206
